$d = $word.ActiveDocument

# Locate the last paragraph in the document (the "...irrespective of amount."
# bullet) and append a brand-new sibling bullet after it at the same list
# level (ilvl=1, numId=2), matching the surrounding sz=28/szCs=28 formatting.
$lastPara = $d.Paragraphs.Last
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# Insert the full sentence as a single run first so it correctly inherits the
# paragraph's run formatting (sz/szCs), then underline just the requested
# substring afterwards - this keeps every resulting run's size formatting
# intact, matching how Word would split the run on a formatting change.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$fullText = "Children can request a Parent to deposit money into the wallet if there isn" + [char]0x2019 + "t enough money to withdraw for their purposes."
$newRange.InsertAfter($fullText)

$underlineRange = $d.Paragraphs.Last.Range.Duplicate
$found = $underlineRange.Find.Execute("request a Parent to deposit money", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$underlineRange.Font.Underline = 1
